$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (ASSISTS) rows 2-41 were stored as inline strings ("0", "1", ...).
# Rewrite them as real numbers, preserving the same values.
for ($r = 2; $r -le 41; $r++) {
    $cell = $ws.Cells.Item($r, 6)  # column F
    $cell.Value = [double]$cell.Text
}

# Column H (CHAMPION): rename "Pyke" to "Evelynn" for the affected rows.
$pykeRows = @(5, 11, 17, 23, 29, 35, 41)
foreach ($r in $pykeRows) {
    $cell = $ws.Cells.Item($r, 8)  # column H
    if ($cell.Text -eq "Pyke") {
        $cell.Value = "Evelynn"
    }
}
